$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 54.69462833333333
$ws.Cells.Item(2, 8).Value = 164.083885
$ws.Cells.Item(2, 9).Value = 0.2790924419198448
$ws.Cells.Item(2, 10).Value = 0.2790924419198448
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.300108
$ws.Cells.Item(2, 14).Value = 0.900324
$ws.Cells.Item(2, 15).Value = 0.04336467165755702
$ws.Cells.Item(2, 16).Value = 0.04336467165755702
$ws.Cells.Item(2, 17).Value = 16.41429551986
$ws.Cells.Item(2, 18).Value = 147.72865967874
$ws.Cells.Item(2, 19).Value = 0.01210275210595987
$ws.Cells.Item(2, 20).Value = 0.01210275210595987

# Row 3
$ws.Cells.Item(3, 7).Value = 54.69462833333333
$ws.Cells.Item(3, 8).Value = 164.083885
$ws.Cells.Item(3, 9).Value = 0.2790924419198448
$ws.Cells.Item(3, 10).Value = 0.2790924419198448
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.310677666666666
$ws.Cells.Item(3, 14).Value = 6.932033
$ws.Cells.Item(3, 15).Value = 0.3338857288757714
$ws.Cells.Item(3, 16).Value = 0.3338857288757713
$ws.Cells.Item(3, 17).Value = 126.3816561764672
$ws.Cells.Item(3, 18).Value = 1137.434905588205
$ws.Cells.Item(3, 19).Value = 0.09318498339412627
$ws.Cells.Item(3, 20).Value = 0.09318498339412626

# Row 4
$ws.Cells.Item(4, 7).Value = 54.69462833333333
$ws.Cells.Item(4, 8).Value = 164.083885
$ws.Cells.Item(4, 9).Value = 0.2790924419198448
$ws.Cells.Item(4, 10).Value = 0.2790924419198448
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.274320333333333
$ws.Cells.Item(4, 14).Value = 12.822961
$ws.Cells.Item(4, 15).Value = 0.617625980694349
$ws.Cells.Item(4, 16).Value = 0.6176259806943489
$ws.Cells.Item(4, 17).Value = 233.7823620092761
$ws.Cells.Item(4, 18).Value = 2104.041258083485
$ws.Cells.Item(4, 19).Value = 0.1723747431451248
$ws.Cells.Item(4, 20).Value = 0.1723747431451247

# Row 5
$ws.Cells.Item(5, 7).Value = 54.69462833333333
$ws.Cells.Item(5, 8).Value = 164.083885
$ws.Cells.Item(5, 9).Value = 0.2790924419198448
$ws.Cells.Item(5, 10).Value = 0.2790924419198448
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.03545833333333333
$ws.Cells.Item(5, 14).Value = 0.106375
$ws.Cells.Item(5, 15).Value = 0.005123618772322663
$ws.Cells.Item(5, 16).Value = 0.005123618772322661
$ws.Cells.Item(5, 17).Value = 1.939380362986111
$ws.Cells.Item(5, 18).Value = 17.454423266875
$ws.Cells.Item(5, 19).Value = 0.001429963274633889
$ws.Cells.Item(5, 20).Value = 0.001429963274633889

# Row 6
$ws.Cells.Item(6, 7).Value = 19.32115333333334
$ws.Cells.Item(6, 8).Value = 57.96346000000001
$ws.Cells.Item(6, 9).Value = 0.09859081282432611
$ws.Cells.Item(6, 10).Value = 0.09859081282432611
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.300108
$ws.Cells.Item(6, 14).Value = 0.900324
$ws.Cells.Item(6, 15).Value = 0.04336467165755702
$ws.Cells.Item(6, 16).Value = 0.04336467165755702
$ws.Cells.Item(6, 17).Value = 5.798432684560002
$ws.Cells.Item(6, 18).Value = 52.18589416104001
$ws.Cells.Item(6, 19).Value = 0.004275358226578564
$ws.Cells.Item(6, 20).Value = 0.004275358226578563

# Row 7
$ws.Cells.Item(7, 7).Value = 19.32115333333334
$ws.Cells.Item(7, 8).Value = 57.96346000000001
$ws.Cells.Item(7, 9).Value = 0.09859081282432611
$ws.Cells.Item(7, 10).Value = 0.09859081282432611
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.310677666666666
$ws.Cells.Item(7, 14).Value = 6.932033
$ws.Cells.Item(7, 15).Value = 0.3338857288757714
$ws.Cells.Item(7, 16).Value = 0.3338857288757713
$ws.Cells.Item(7, 17).Value = 44.64495750157556
$ws.Cells.Item(7, 18).Value = 401.8046175141801
$ws.Cells.Item(7, 19).Value = 0.03291806540030487
$ws.Cells.Item(7, 20).Value = 0.03291806540030487

# Row 8
$ws.Cells.Item(8, 7).Value = 19.32115333333334
$ws.Cells.Item(8, 8).Value = 57.96346000000001
$ws.Cells.Item(8, 9).Value = 0.09859081282432611
$ws.Cells.Item(8, 10).Value = 0.09859081282432611
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 4.274320333333333
$ws.Cells.Item(8, 14).Value = 12.822961
$ws.Cells.Item(8, 15).Value = 0.617625980694349
$ws.Cells.Item(8, 16).Value = 0.6176259806943489
$ws.Cells.Item(8, 17).Value = 82.58479855611779
$ws.Cells.Item(8, 18).Value = 743.2631870050601
$ws.Cells.Item(8, 19).Value = 0.06089224745807741
$ws.Cells.Item(8, 20).Value = 0.06089224745807741

# Row 9
$ws.Cells.Item(9, 7).Value = 19.32115333333334
$ws.Cells.Item(9, 8).Value = 57.96346000000001
$ws.Cells.Item(9, 9).Value = 0.09859081282432611
$ws.Cells.Item(9, 10).Value = 0.09859081282432611
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.03545833333333333
$ws.Cells.Item(9, 14).Value = 0.106375
$ws.Cells.Item(9, 15).Value = 0.005123618772322663
$ws.Cells.Item(9, 16).Value = 0.005123618772322661
$ws.Cells.Item(9, 17).Value = 0.685095895277778
$ws.Cells.Item(9, 18).Value = 6.165863057500001
$ws.Cells.Item(9, 19).Value = 0.0005051417393652672
$ws.Cells.Item(9, 20).Value = 0.000505141739365267

# Row 10
$ws.Cells.Item(10, 7).Value = 11.023718
$ws.Cells.Item(10, 8).Value = 33.071154
$ws.Cells.Item(10, 9).Value = 0.05625116157486912
$ws.Cells.Item(10, 10).Value = 0.05625116157486911
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.300108
$ws.Cells.Item(10, 14).Value = 0.900324
$ws.Cells.Item(10, 15).Value = 0.04336467165755702
$ws.Cells.Item(10, 16).Value = 0.04336467165755702
$ws.Cells.Item(10, 17).Value = 3.308305961544
$ws.Cells.Item(10, 18).Value = 29.774753653896
$ws.Cells.Item(10, 19).Value = 0.002439313152050388
$ws.Cells.Item(10, 20).Value = 0.002439313152050387

# Row 11
$ws.Cells.Item(11, 7).Value = 11.023718
$ws.Cells.Item(11, 8).Value = 33.071154
$ws.Cells.Item(11, 9).Value = 0.05625116157486912
$ws.Cells.Item(11, 10).Value = 0.05625116157486911
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 2.310677666666666
$ws.Cells.Item(11, 14).Value = 6.932033
$ws.Cells.Item(11, 15).Value = 0.3338857288757714
$ws.Cells.Item(11, 16).Value = 0.3338857288757713
$ws.Cells.Item(11, 17).Value = 25.47225898623133
$ws.Cells.Item(11, 18).Value = 229.250330876082
$ws.Cells.Item(11, 19).Value = 0.01878146008253396
$ws.Cells.Item(11, 20).Value = 0.01878146008253396

# Row 12
$ws.Cells.Item(12, 7).Value = 11.023718
$ws.Cells.Item(12, 8).Value = 33.071154
$ws.Cells.Item(12, 9).Value = 0.05625116157486912
$ws.Cells.Item(12, 10).Value = 0.05625116157486911
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 4.274320333333333
$ws.Cells.Item(12, 14).Value = 12.822961
$ws.Cells.Item(12, 15).Value = 0.617625980694349
$ws.Cells.Item(12, 16).Value = 0.6176259806943489
$ws.Cells.Item(12, 17).Value = 47.11890199633266
$ws.Cells.Item(12, 18).Value = 424.070117966994
$ws.Cells.Item(12, 19).Value = 0.03474217883287482
$ws.Cells.Item(12, 20).Value = 0.03474217883287481

# Row 13
$ws.Cells.Item(13, 7).Value = 11.023718
$ws.Cells.Item(13, 8).Value = 33.071154
$ws.Cells.Item(13, 9).Value = 0.05625116157486912
$ws.Cells.Item(13, 10).Value = 0.05625116157486911
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.03545833333333333
$ws.Cells.Item(13, 14).Value = 0.106375
$ws.Cells.Item(13, 15).Value = 0.005123618772322663
$ws.Cells.Item(13, 16).Value = 0.005123618772322661
$ws.Cells.Item(13, 17).Value = 0.3908826674166667
$ws.Cells.Item(13, 18).Value = 3.51794400675
$ws.Cells.Item(13, 19).Value = 0.0002882095074099547
$ws.Cells.Item(13, 20).Value = 0.0002882095074099545

# Row 14
$ws.Cells.Item(14, 7).Value = 110.9336623333333
$ws.Cells.Item(14, 8).Value = 332.800987
$ws.Cells.Item(14, 9).Value = 0.5660655836809599
$ws.Cells.Item(14, 10).Value = 0.5660655836809599
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.300108
$ws.Cells.Item(14, 14).Value = 0.900324
$ws.Cells.Item(14, 15).Value = 0.04336467165755702
$ws.Cells.Item(14, 16).Value = 0.04336467165755702
$ws.Cells.Item(14, 17).Value = 33.292079535532
$ws.Cells.Item(14, 18).Value = 299.628715819788
$ws.Cells.Item(14, 19).Value = 0.0245472481729682
$ws.Cells.Item(14, 20).Value = 0.02454724817296819

# Row 15
$ws.Cells.Item(15, 7).Value = 110.9336623333333
$ws.Cells.Item(15, 8).Value = 332.800987
$ws.Cells.Item(15, 9).Value = 0.5660655836809599
$ws.Cells.Item(15, 10).Value = 0.5660655836809599
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 2.310677666666666
$ws.Cells.Item(15, 14).Value = 6.932033
$ws.Cells.Item(15, 15).Value = 0.3338857288757714
$ws.Cells.Item(15, 16).Value = 0.3338857288757713
$ws.Cells.Item(15, 17).Value = 256.3319360351745
$ws.Cells.Item(15, 18).Value = 2306.98742431657
$ws.Cells.Item(15, 19).Value = 0.1890012199988063
$ws.Cells.Item(15, 20).Value = 0.1890012199988062

# Row 16
$ws.Cells.Item(16, 7).Value = 110.9336623333333
$ws.Cells.Item(16, 8).Value = 332.800987
$ws.Cells.Item(16, 9).Value = 0.5660655836809599
$ws.Cells.Item(16, 10).Value = 0.5660655836809599
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 4.274320333333333
$ws.Cells.Item(16, 14).Value = 12.822961
$ws.Cells.Item(16, 15).Value = 0.617625980694349
$ws.Cells.Item(16, 16).Value = 0.6176259806943489
$ws.Cells.Item(16, 17).Value = 474.1660085625007
$ws.Cells.Item(16, 18).Value = 4267.494077062506
$ws.Cells.Item(16, 19).Value = 0.349616811258272
$ws.Cells.Item(16, 20).Value = 0.3496168112582719

# Row 17
$ws.Cells.Item(17, 7).Value = 110.9336623333333
$ws.Cells.Item(17, 8).Value = 332.800987
$ws.Cells.Item(17, 9).Value = 0.5660655836809599
$ws.Cells.Item(17, 10).Value = 0.5660655836809599
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.03545833333333333
$ws.Cells.Item(17, 14).Value = 0.106375
$ws.Cells.Item(17, 15).Value = 0.005123618772322663
$ws.Cells.Item(17, 16).Value = 0.005123618772322661
$ws.Cells.Item(17, 17).Value = 3.933522776902777
$ws.Cells.Item(17, 18).Value = 35.401704992125
$ws.Cells.Item(17, 19).Value = 0.002900304250913552
$ws.Cells.Item(17, 20).Value = 0.002900304250913551
